# Updates the "Team of Outs" roster sheet:
#  - Removes the "Keon Ellis" row (player no longer listed)
#  - Reorders the remaining player rows into the new sequence

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Locate and remove the row for Keon Ellis. Note: this runtime's bare
# `.Value` getter is unreliable for reads (it can return a reflection
# description instead of the cell contents) - use `.Value2` when reading.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($r = $lastRow; $r -ge 2; $r--) {
    if ($ws.Cells.Item($r, 1).Value2 -eq "Keon Ellis") {
        $ws.Rows.Item($r).Delete()
    }
}

# Final desired order of the remaining 16 players (name, position, team)
$data = @(
    @("Ja Morant",        "PG",        "Memphis Grizzlies"),
    @("Tyler Herro",      "PG,SG",     "Miami Heat"),
    @("Mikal Bridges",    "SG,SF,PF",  "New York Knicks"),
    @("T.J. McConnell",   "PG",        "Indiana Pacers"),
    @("Scottie Barnes",   "SG,SF,PF",  "Toronto Raptors"),
    @("P.J. Washington",  "PF",        "Dallas Mavericks"),
    @("Santi Aldama",     "PF,C",      "Memphis Grizzlies"),
    @("Miles Bridges",    "SF,PF",     "Charlotte Hornets"),
    @("Brook Lopez",      "C",         "Milwaukee Bucks"),
    @("Evan Mobley",      "PF,C",      "Cleveland Cavaliers"),
    @("Andrew Nembhard",  "PG,SG",     "Indiana Pacers"),
    @("Nikola Vucevic",   "PF,C",      "Chicago Bulls"),
    @("De'Aaron Fox",     "PG",        "Sacramento Kings"),
    @("DeMar DeRozan",    "SF,PF",     "Sacramento Kings"),
    @("Luka Doncic",      "PG,SG",     "Dallas Mavericks"),
    @("Josh Giddey",      "PG,SG,SF",  "Chicago Bulls")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
